$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a brand-new row above row 13; this shifts rows 13..27 down to 14..28,
# Excel auto-adjusts all the relative formulas in the shifted rows, and the new
# row 13 inherits its formatting from row 12 above it (same "data row" layout).
$ws.Rows.Item(13).Insert()

# New experimental data row.
$ws.Range("C13").Value = 0.2953
$ws.Range("D13").Value = 0.1171
$ws.Range("E13").Value = 0.5876

$ws.Range("G13").Formula = "=C13*`$A`$2/(C13*`$A`$2+E13*`$A`$4+D13*`$A`$6)"
$ws.Range("H13").Formula = "=D13*`$A`$6/(C13*`$A`$2+E13*`$A`$4+D13*`$A`$6)"
$ws.Range("I13").Formula = "=E13*`$A`$4/(C13*`$A`$2+E13*`$A`$4+D13*`$A`$6)"
$ws.Range("K13").Formula = "=SUM(G13:I13)"
$ws.Range("L13").Formula = "=G13*I13"

# K/L stay on the workbook's default (unstyled) format, same as K12/L12 --
# restore that *after* the formula writes above, which otherwise pick up I13's
# number format as they're typed next to it.
$ws.Range("K13:L13").Style = "Normal"

# New column N (14) sized to fit its new best-fit content.
$ws.Columns.Item(14).ColumnWidth = 17.85546875

# Two brand-new summary rows at the bottom of the sheet.
$ws.Range("N34").Formula = "=D15+E15"
$ws.Range("N34").NumberFormat = "0.000000000000000"

$ws.Range("N37").Formula = "=5876+1171"

$ws.Range("A6").Select()

$wb.Save()
